$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (context G=5533)
$ws.Range("H11").Value = 100009.3
$ws.Range("I11").Value = 100009.3
$ws.Range("K11").Value = 100009.3
$ws.Range("M11").Value = -99869.3

# Row 17 (context G=38956)
$ws.Range("H17").Value = 1044.4546
$ws.Range("J17").Value = 1044.4546
$ws.Range("L17").Value = 3133.3638
$ws.Range("N17").Value = -3469.3638

# Row 43 (context G=5472)
$ws.Range("H43").Value = 2366.6667
$ws.Range("J43").Value = 799.75
$ws.Range("L43").Value = 799.75
$ws.Range("N43").Value = -937.75

# Row 62 (context G=27781)
$ws.Range("H62").Value = 2736
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3226.6667
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 3226.6667
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4474.6667

# Row 65 (context G=27781)
$ws.Range("H65").Value = 2736
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3226.6667
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 16133.3335
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -22373.3335

# Row 98 (context G=36237)
$ws.Range("H98").Value = 800.2857
$ws.Range("I98").Value = 800.2857
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 800.2857
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 697.7143
$ws.Range("N98").ClearContents()

# Row 112 (context G=27960)
$ws.Range("H112").Value = 1796.1333
$ws.Range("J112").Value = 1887.8462
$ws.Range("L112").Value = 5663.5386
$ws.Range("N112").Value = -7879.5386

# Row 122 (context G=36237)
$ws.Range("H122").Value = 800.2857
$ws.Range("I122").Value = 800.2857
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2400.8571
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 49.14289999999983
$ws.Range("N122").ClearContents()

# Row 135 (context G=44047)
$ws.Range("H135").Value = 1868.875
$ws.Range("I135").Value = 1045.8
$ws.Range("J135").Value = 2456.7856
$ws.Range("K135").Value = 9412.199999999999
$ws.Range("L135").Value = 22111.0704
$ws.Range("M135").Value = -6877.199999999999
$ws.Range("N135").Value = -27181.0704

# Row 137 (context G=44013)
$ws.Range("H137").Value = 1664.625
$ws.Range("I137").Value = 1668.8422
$ws.Range("J137").Value = 1658.4615
$ws.Range("K137").Value = 5006.5266
$ws.Range("L137").Value = 4975.3845
$ws.Range("M137").Value = -2456.5266
$ws.Range("N137").Value = -10075.3845


$ws = $wb.Worksheets.Item("ARM")
# Row 2 (context G=27713)
$ws.Range("H2").Value = 144699.14
$ws.Range("I2").Value = 2620
$ws.Range("K2").Value = 2620
$ws.Range("M2").Value = -2507

# Row 32 (context G=44147)
$ws.Range("H32").Value = 25659.303
$ws.Range("I32").Value = 5733.7354
$ws.Range("J32").Value = 195026.62
$ws.Range("K32").Value = 5733.7354
$ws.Range("L32").Value = 195026.62
$ws.Range("M32").Value = -5446.7354
$ws.Range("N32").Value = -195600.62

# Row 45 (context G=27714)
$ws.Range("H45").Value = 72246.64
$ws.Range("I45").Value = 167470.17
$ws.Range("K45").Value = 167470.17
$ws.Range("M45").Value = -167093.17

# Row 116 (context G=27713)
$ws.Range("H116").Value = 144699.14
$ws.Range("I116").Value = 2620
$ws.Range("K116").Value = 2620
$ws.Range("M116").Value = -326

# Row 122 (context G=36168)
$ws.Range("H122").Value = 1271.2646
$ws.Range("I122").Value = 1228.5769
$ws.Range("K122").Value = 3685.7307
$ws.Range("M122").Value = -1235.7307

# Row 132 (context G=43997)
$ws.Range("H132").Value = 4071.6
$ws.Range("I132").Value = 4200.476
$ws.Range("J132").Value = 3770.889
$ws.Range("K132").Value = 12601.428
$ws.Range("L132").Value = 11312.667
$ws.Range("M132").Value = -10071.428
$ws.Range("N132").Value = -16372.667


$ws = $wb.Worksheets.Item("BSM")
# Row 3 (context G=27713)
$ws.Range("H3").Value = 144699.14
$ws.Range("I3").Value = 2620
$ws.Range("K3").Value = 2620
$ws.Range("M3").Value = -2506

# Row 94 (context G=19939)
$ws.Range("H94").Value = 454.13333
$ws.Range("I94").Value = 411.05264
$ws.Range("K94").Value = 411.05264
$ws.Range("M94").Value = 39.94736


$ws = $wb.Worksheets.Item("CRP")
# Row 31 (context G=44023)
$ws.Range("H31").Value = 20721.773
$ws.Range("I31").Value = 1428.68
$ws.Range("J31").Value = 30368.32
$ws.Range("K31").Value = 1428.68
$ws.Range("L31").Value = 30368.32
$ws.Range("M31").Value = -1133.68
$ws.Range("N31").Value = -30958.32

# Row 34 (context G=44023)
$ws.Range("H34").Value = 20721.773
$ws.Range("I34").Value = 1428.68
$ws.Range("J34").Value = 30368.32
$ws.Range("K34").Value = 1428.68
$ws.Range("L34").Value = 30368.32
$ws.Range("M34").Value = -1226.68
$ws.Range("N34").Value = -30772.32

# Row 94 (context G=32934)
$ws.Range("H94").Value = 1107.375
$ws.Range("I94").Value = 965
$ws.Range("J94").Value = 1154.8334
$ws.Range("K94").Value = 965
$ws.Range("L94").Value = 1154.8334
$ws.Range("M94").Value = -514
$ws.Range("N94").Value = -2056.8334

# Row 122 (context G=36196)
$ws.Range("H122").Value = 387.10526
$ws.Range("I122").Value = 207.875
$ws.Range("J122").Value = 517.4545000000001
$ws.Range("K122").Value = 623.625
$ws.Range("L122").Value = 1552.3635
$ws.Range("M122").Value = 1826.375
$ws.Range("N122").Value = -6452.3635

# Row 132 (context G=44019)
$ws.Range("H132").Value = 125007496
$ws.Range("I132").Value = 250014300
$ws.Range("K132").Value = 750042900
$ws.Range("M132").Value = -750040370


$ws = $wb.Worksheets.Item("CUL")
# Row 123 (context G=36037)
$ws.Range("H123").Value = 5700
$ws.Range("J123").Value = 5700
$ws.Range("L123").Value = 17100
$ws.Range("N123").Value = -22000


$ws = $wb.Worksheets.Item("GSM")
# Row 102 (context G=36169)
$ws.Range("H102").Value = 2715.8235
$ws.Range("I102").Value = 1808.091
$ws.Range("K102").Value = 1808.091
$ws.Range("M102").Value = -186.0909999999999

# Row 122 (context G=36182)
$ws.Range("H122").Value = 946.0357
$ws.Range("I122").Value = 872.8570999999999
$ws.Range("J122").Value = 1019.2143
$ws.Range("K122").Value = 2618.5713
$ws.Range("L122").Value = 3057.6429
$ws.Range("M122").Value = -168.5712999999996
$ws.Range("N122").Value = -7957.6429


$ws = $wb.Worksheets.Item("LTW")
# Row 40 (context G=36248)
$ws.Range("H40").Value = 47245.637
$ws.Range("I40").Value = 101165.7
$ws.Range("K40").Value = 101165.7
$ws.Range("M40").Value = -101029.7

# Row 61 (context G=27740)
$ws.Range("H61").Value = 1746.8
$ws.Range("I61").Value = 1607.5385
$ws.Range("K61").Value = 1607.5385
$ws.Range("M61").Value = -1405.5385

# Row 113 (context G=27740)
$ws.Range("H113").Value = 1746.8
$ws.Range("I113").Value = 1607.5385
$ws.Range("K113").Value = 1607.5385
$ws.Range("M113").Value = 562.4614999999999

# Row 132 (context G=44058)
$ws.Range("H132").Value = 4391.375
$ws.Range("I132").Value = 4468.1577
$ws.Range("J132").Value = 4099.6
$ws.Range("K132").Value = 13404.4731
$ws.Range("L132").Value = 12298.8
$ws.Range("M132").Value = -10874.4731
$ws.Range("N132").Value = -17358.8


$ws = $wb.Worksheets.Item("WVR")
# Row 113 (context G=27752)
$ws.Range("H113").Value = 848.9375
$ws.Range("I113").Value = 685.4286
$ws.Range("J113").Value = 976.1111
$ws.Range("K113").Value = 2056.2858
$ws.Range("L113").Value = 2928.3333
$ws.Range("M113").Value = 113.7142000000003
$ws.Range("N113").Value = -7268.3333

# Row 122 (context G=36208)
$ws.Range("H122").Value = 1309
$ws.Range("I122").Value = 1188.1
$ws.Range("J122").Value = 1577.6666
$ws.Range("K122").Value = 3564.3
$ws.Range("L122").Value = 4732.9998
$ws.Range("M122").Value = -1114.3
$ws.Range("N122").Value = -9632.9998

# Row 132 (context G=44029)
$ws.Range("H132").Value = 2337.0256
$ws.Range("I132").Value = 2371.3547
$ws.Range("K132").Value = 7114.0641
$ws.Range("M132").Value = -4584.0641

